# [Fonds de solidarite] Add 2020-12-16 data
# The source data cells in columns C/D/E are stored as text (not numbers),
# so we force the Text number format before writing each value to avoid
# Excel auto-converting the numeric-looking strings into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Row 9 - Bourgogne-Franche-Comte / Entrepreneur individuel
Set-TextValue "C9" "320"
Set-TextValue "D9" "300"
Set-TextValue "E9" "1088667.48"

# Row 11 - Bourgogne-Franche-Comte / SARL
Set-TextValue "C11" "551"
Set-TextValue "D11" "491"
Set-TextValue "E11" "4231204.00"

# Row 12 - Bourgogne-Franche-Comte / SAS
Set-TextValue "C12" "275"
Set-TextValue "E12" "2184224.49"

# Row 34 - Grand Est / SARL
Set-TextValue "C34" "898"
Set-TextValue "E34" "7313713.66"

# Row 51 - Hauts-de-France / SARL
Set-TextValue "C51" "1181"
Set-TextValue "E51" "9030924.05"

# Row 52 - Hauts-de-France / SAS
Set-TextValue "C52" "814"
Set-TextValue "E52" "5619450.79"

# Row 60 - Ile-de-France / SAS
Set-TextValue "C60" "6749"
Set-TextValue "E60" "31007350.12"

# Row 65 - La Reunion / SAS
Set-TextValue "C65" "60"
Set-TextValue "D65" "60"
Set-TextValue "E65" "425050.00"
